# Update cryptocurrency price (D) and 1h volume change (E) columns
# to reflect the latest scraped values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.703.56"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.678.09"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3919"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3968"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.005"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "51.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.409"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08628"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.326"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.785"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001318"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.04%  "
$ws.Range("D17").Value = "1.689.28"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07080"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.079"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").Value = "24.708.08"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.362"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "23.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.764"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.778"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "149.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.845"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.402"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.34%  "
$ws.Range("D33").Value = "1.881.39"
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08421"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03074"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.919"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2788"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09464"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7917"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.44%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7127"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.560"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.172"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08671"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.336"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.69%  "
